$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 20000
$ws.Range("I10").Value = 20000
$ws.Range("K10").Value = 20000
$ws.Range("M10").Value = -19707
$ws.Range("H129").Value = 1138.881
$ws.Range("J129").Value = 1228.2433
$ws.Range("L129").Value = 3684.7299
$ws.Range("N129").Value = -13684.7299
$ws.Range("H137").Value = 1184.1025
$ws.Range("I137").Value = 1140
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 3420
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -870
$ws.Range("N137").Value = -11100

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1815
$ws.Range("I61").Value = 1576.7307
$ws.Range("J61").Value = 2700
$ws.Range("K61").Value = 1576.7307
$ws.Range("L61").Value = 2700
$ws.Range("M61").Value = -1364.7307
$ws.Range("N61").Value = -3124
$ws.Range("H74").Value = 1316.5454
$ws.Range("I74").Value = 926.7143
$ws.Range("J74").Value = 1998.75
$ws.Range("K74").Value = 926.7143
$ws.Range("L74").Value = 1998.75
$ws.Range("M74").Value = -52.71429999999998
$ws.Range("N74").Value = -3746.75
$ws.Range("H77").Value = 1316.5454
$ws.Range("I77").Value = 926.7143
$ws.Range("J77").Value = 1998.75
$ws.Range("K77").Value = 4633.5715
$ws.Range("L77").Value = 9993.75
$ws.Range("M77").Value = -265.5715
$ws.Range("N77").Value = -18729.75
$ws.Range("H131").Value = 38500.5
$ws.Range("J131").Value = 38500.5
$ws.Range("L131").Value = 38500.5
$ws.Range("N131").Value = -48580.5
$ws.Range("H132").Value = 4007.196
$ws.Range("I132").Value = 4781.9375
$ws.Range("J132").Value = 2702.3684
$ws.Range("K132").Value = 14345.8125
$ws.Range("L132").Value = 8107.1052
$ws.Range("M132").Value = -11815.8125
$ws.Range("N132").Value = -13167.1052
$ws.Range("H135").Value = 53598
$ws.Range("J135").Value = 53598
$ws.Range("L135").Value = 53598
$ws.Range("N135").Value = -63738
$ws.Range("H136").Value = 1815
$ws.Range("I136").Value = 1576.7307
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 4730.1921
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -2180.1921
$ws.Range("N136").Value = -13200

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2469306.5
$ws.Range("I80").Value = 9259292
$ws.Range("J80").Value = 220.81818
$ws.Range("K80").Value = 9259292
$ws.Range("L80").Value = 220.81818
$ws.Range("M80").Value = -9258294
$ws.Range("N80").Value = -2216.81818
$ws.Range("H83").Value = 2469306.5
$ws.Range("I83").Value = 9259292
$ws.Range("J83").Value = 220.81818
$ws.Range("K83").Value = 46296460
$ws.Range("L83").Value = 1104.0909
$ws.Range("M83").Value = -46291468
$ws.Range("N83").Value = -11088.0909
$ws.Range("H134").Value = 2144.0435
$ws.Range("I134").Value = 1735.1177
$ws.Range("K134").Value = 5205.3531
$ws.Range("M134").Value = -2670.3531

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 17029.334
$ws.Range("I12").Value = 644
$ws.Range("J12").Value = 49800
$ws.Range("K12").Value = 644
$ws.Range("L12").Value = 49800
$ws.Range("M12").Value = -474
$ws.Range("N12").Value = -50140
$ws.Range("H31").Value = 2059.0454
$ws.Range("I31").Value = 1594.6842
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1594.6842
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1299.6842
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 2059.0454
$ws.Range("I34").Value = 1594.6842
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1594.6842
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -1392.6842
$ws.Range("N34").Value = -5404
$ws.Range("H58").Value = 863193.75
$ws.Range("I58").Value = 1324447.4
$ws.Range("J58").Value = 2186.9333
$ws.Range("K58").Value = 1324447.4
$ws.Range("L58").Value = 2186.9333
$ws.Range("M58").Value = -1324244.4
$ws.Range("N58").Value = -2592.9333
$ws.Range("H132").Value = 411453.9
$ws.Range("I132").Value = 484042.2
$ws.Range("J132").Value = 4959.6
$ws.Range("K132").Value = 1452126.6
$ws.Range("L132").Value = 14878.8
$ws.Range("M132").Value = -1449596.6
$ws.Range("N132").Value = -19938.8
$ws.Range("H134").Value = 2352.5833
$ws.Range("I134").Value = 1667.8334
$ws.Range("K134").Value = 5003.5002
$ws.Range("M134").Value = -2468.5002
$ws.Range("H136").Value = 863193.75
$ws.Range("I136").Value = 1324447.4
$ws.Range("J136").Value = 2186.9333
$ws.Range("K136").Value = 3973342.2
$ws.Range("L136").Value = 6560.7999
$ws.Range("M136").Value = -3970792.2
$ws.Range("N136").Value = -11660.7999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1905.9286
$ws.Range("J5").Value = 833
$ws.Range("L5").Value = 2499
$ws.Range("N5").Value = -2723
$ws.Range("H104").Value = 6587.875
$ws.Range("J104").Value = 6587.875
$ws.Range("L104").Value = 19763.625
$ws.Range("N104").Value = -25005.625
$ws.Range("H135").Value = 1905.9286
$ws.Range("J135").Value = 833
$ws.Range("L135").Value = 7497
$ws.Range("N135").Value = -12567

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 33466666
$ws.Range("I7").Value = 100000000
$ws.Range("K7").Value = 100000000
$ws.Range("M7").Value = -99999888
$ws.Range("H8").Value = 33466666
$ws.Range("I8").Value = 100000000
$ws.Range("K8").Value = 100000000
$ws.Range("M8").Value = -99999861
$ws.Range("H109").Value = 9020.823
$ws.Range("J109").Value = 9020.823
$ws.Range("L109").Value = 9020.823
$ws.Range("N109").Value = -11100.823
$ws.Range("H131").Value = 41327.25
$ws.Range("J131").Value = 41327.25
$ws.Range("L131").Value = 41327.25
$ws.Range("N131").Value = -51407.25
$ws.Range("H132").Value = 1806.5918
$ws.Range("I132").Value = 1203.5
$ws.Range("K132").Value = 3610.5
$ws.Range("M132").Value = -1080.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 20033334
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H131").Value = 30324
$ws.Range("J131").Value = 30324
$ws.Range("L131").Value = 30324
$ws.Range("N131").Value = -40404
$ws.Range("H132").Value = 4342.36
$ws.Range("I132").Value = 4713.6
$ws.Range("J132").Value = 3785.5
$ws.Range("K132").Value = 14140.8
$ws.Range("L132").Value = 11356.5
$ws.Range("M132").Value = -11610.8
$ws.Range("N132").Value = -16416.5
$ws.Range("H136").Value = 29708132
$ws.Range("I136").Value = 40001500
$ws.Range("J136").Value = 1115445
$ws.Range("K136").Value = 120004500
$ws.Range("L136").Value = 3346335
$ws.Range("M136").Value = -120001950
$ws.Range("N136").Value = -3351435

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496
$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716
$ws.Range("H123").Value = 23299.39
$ws.Range("J123").Value = 23299.39
$ws.Range("L123").Value = 23299.39
$ws.Range("N123").Value = -33099.39
$ws.Range("H132").Value = 3189.5264
$ws.Range("I132").Value = 2644
$ws.Range("J132").Value = 3507.75
$ws.Range("K132").Value = 7932
$ws.Range("L132").Value = 10523.25
$ws.Range("M132").Value = -5402
$ws.Range("N132").Value = -15583.25
$ws.Range("H136").Value = 2132.52
$ws.Range("I136").Value = 1979.6316
$ws.Range("J136").Value = 2616.6667
$ws.Range("K136").Value = 5938.8948
$ws.Range("L136").Value = 7850.000100000001
$ws.Range("M136").Value = -3388.8948
$ws.Range("N136").Value = -12950.0001
